$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.022.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.619.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.45%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.620.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.000.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0740"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.44%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("E21").Value = "  -0.74%  "

$ws.Range("E22").Value = "  -5.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.73%  "

$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.336.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.94%  "

$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("E35").Value = "  -0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0175"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.51%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.98%  "

$ws.Range("E43").Value = "  -2.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.755.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.842"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +26.69%  "

$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0993"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.33%  "
